$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.381.10"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "3.503.61"
$ws.Range("E3").Value = "  -2.67%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'606.09"
$ws.Range("E5").Value = "  -2.74%  "
$ws.Range("D6").Value = "'149.20"
$ws.Range("E6").Value = "  -4.41%  "
$ws.Range("D7").Value = "3.502.92"
$ws.Range("E7").Value = "  -2.60%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").Value = "'7.04"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("D14").Value = "4.096.25"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").Value = "3.503.01"
$ws.Range("E16").Value = "  -3.54%  "
$ws.Range("D17").Value = "67.369.62"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").Value = "'15.15"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("D21").Value = "'445.54"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").Value = "'9.18"
$ws.Range("E22").Value = "  -6.54%  "
$ws.Range("D23").Value = "'0.623"
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("D24").Value = "'77.37"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "3.643.62"
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'0.0000127"
$ws.Range("E27").Value = "  +8.52%  "
$ws.Range("D28").Value = "'10.31"
$ws.Range("E28").Value = "  -4.09%  "
$ws.Range("D29").Value = "'8.34"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").Value = "'2.49"
$ws.Range("E30").Value = "  -4.00%  "
$ws.Range("E32").Value = "  -7.20%  "
$ws.Range("D33").Value = "'0.164"
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("D34").Value = "'25.67"
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").Value = "3.495.49"
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("E37").Value = "  -4.85%  "
$ws.Range("D38").Value = "'8.04"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'177.42"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "'0.0869"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  -3.96%  "
$ws.Range("D45").Value = "'0.877"
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").Value = "'27.30"
$ws.Range("E47").Value = "  -3.69%  "
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").Value = "'2.55"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").Value = "'0.989"
$ws.Range("E51").Value = "  -2.24%  "
